# class MonthChecker added, which check which month is to be reported.
# ReportWriter now writes only the forecasted Turn Over values for the
# desired month into column D (rows 4..34).
#
# Each value is written as TEXT so that it is stored as a shared string
# (matching the original workbook layout) instead of being auto-converted
# to a number by the COM layer. Using Formula + Copy/PasteSpecial(values)
# bakes the formula result back down to a plain value, so no formula and
# no extra number-format style survive in the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$text)
    $cell.Formula = '=TEXT("' + $text + '","@")'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$forecast = @{
    4  = "25832"
    5  = "37363"
    6  = "22131"
    7  = "150000.0"
    8  = "130000.0"
    9  = "120000.0"
    10 = "120000.0"
    11 = "130000.0"
    12 = "150000.0"
    13 = "26398"
    14 = "150000.0"
    15 = "130000.0"
    16 = "120000.0"
    17 = "120000.0"
    18 = "130000.0"
    19 = "150000.0"
    20 = "22447"
    21 = "121499"
    22 = "123140"
    23 = "142377"
    24 = "133942"
    25 = "147431"
    26 = "221684"
    27 = "25003"
    28 = "148766"
    29 = "120000"
    30 = "100000"
    31 = "140000"
    32 = "250000"
    33 = "300000"
    34 = "14000"
}

for ($row = 4; $row -le 34; $row++) {
    $value = $forecast[$row]
    $cell = $ws.Cells.Item($row, 4)
    Set-TextValue $cell $value
}
